$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Integrated Register")

# --- Row 13: fill in owner/origin/stakeholders/category + probability/impact ---
$ws.Range("G13").Value = "Tim"
$ws.Range("H13").Value = "Team"
$ws.Range("I13").Value = "eCL team"
$ws.Range("J13").Value = "System"
$ws.Range("K13").Value = 0.1
$ws.Range("L13").Value = 2
# Re-apply the Risk Exposure Ranking formula so it recalculates now that
# Probability / Impact are numeric instead of blank (was #VALUE!).
$ws.Range("N13").Formula = "=Table_owssvr_23[[#This Row],[Probability of Occurrence (%) ]]*Table_owssvr_23[[#This Row],[Impact  `nRating   (1-5)]]*5"

# --- Row 14: same fields, different risk owner contact (Jourdain) ---
$ws.Range("G14").Value = "Tim"
$ws.Range("H14").Value = "Jourdain"
$ws.Range("I14").Value = "eCL team"
$ws.Range("J14").Value = "System"
$ws.Range("K14").Value = 0.1
$ws.Range("L14").Value = 2
$ws.Range("N14").Formula = "=Table_owssvr_23[[#This Row],[Probability of Occurrence (%) ]]*Table_owssvr_23[[#This Row],[Impact  `nRating   (1-5)]]*5"

# --- Row 15: new risk description ---
$ws.Range("C15").Value = "Coaching logs are being entered with PII"

# --- Update the saved cursor/selection on the Integrated Register sheet ---
[void]$ws.Range("D15").Select()
